$d = $word.ActiveDocument

$replacements = @(
    @{old="2024-11-06 Wednesday"; new="2024-11-07 Thursday"},
    @{old="823×4=3292"; new="671×8=5368"},
    @{old="856×4=3424"; new="688×2=1376"},
    @{old="451×8=3608"; new="395×4=1580"},
    @{old="167×9=1503"; new="180×7=1260"},
    @{old="186×3=558";  new="191×7=1337"},
    @{old="517×3=1551"; new="103×6=618"},
    @{old="419×3=1257"; new="194×4=776"},
    @{old="546×8=4368"; new="956×9=8604"},
    @{old="841×7=5887"; new="330×2=660"},
    @{old="713×2=1426"; new="703×8=5624"},
    @{old="878×8=7024"; new="981×5=4905"},
    @{old="822×9=7398"; new="441×2=882"},
    @{old="454×3=1362"; new="725×5=3625"},
    @{old="771×5=3855"; new="723×9=6507"},
    @{old="841×4=3364"; new="615×3=1845"},
    @{old="992×6=5952"; new="925×3=2775"},
    @{old="648×8=5184"; new="155×7=1085"},
    @{old="290×7=2030"; new="201×7=1407"},
    @{old="259×5=1295"; new="444×5=2220"},
    @{old="106×9=954";  new="259×9=2331"},
    @{old="665×5=3325"; new="414×6=2484"},
    @{old="909×2=1818"; new="670×6=4020"},
    @{old="528×6=3168"; new="607×7=4249"},
    @{old="517×8=4136"; new="230×4=920"},
    @{old="406×3=1218"; new="568×7=3976"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
